# Updateret tidsregistrering d. 26
# Adds the 26-02-2015 (serial 42061) time-registration entries and tweaks
# the existing "Team 7" note to spell out which use-cases it covered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# --- New row 7: 26/02 - 0.5h - Business-process analyst - Team 5 (kunden) UC5
# Written before the Team-7 edit below so new shared strings are appended
# in the same order the workbook author typed them in.
$ws.Cells.Item(7, 1).Value = 42061
$ws.Cells.Item(7, 2).Value = 0.5
$ws.Cells.Item(7, 3).Value = "Business-process analyst"
$ws.Cells.Item(7, 4).Value = "Samtale med Team 5 (kunden) UC5"

# --- Clarify the existing Team 7 conversation note (row 5, column D)
$ws.Cells.Item(5, 4).Value = "Samtale med Team 7 (kunden) UC1+UC2"

# --- New row 8: 26/02 - 1.2h - System Analyst - domain-model writeup
$ws.Cells.Item(8, 1).Value = 42061
$ws.Cells.Item(8, 2).Value = 1.2
$ws.Cells.Item(8, 3).Value = "System Analyst"
$ws.Cells.Item(8, 4).Value = "Lavet udkast til domænemodel for UC1/UC2 + dataordborg for UC1/UC2. Har gennemgået rettelser med anden gruppe af samme opgaver."

# Leave the selection on the next free row, like the author did after typing.
[void]$ws.Range("A9").Select()
